$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.859.16'
$ws.Range("E2").Value = '  -0.74%  '

# Row 3
$ws.Range("D3").Value = '2.448.20'
$ws.Range("E3").Value = '  +0.38%  '

# Row 4
$ws.Range("E4").Value = '  -0.15%  '

# Row 5
$ws.Range("D5").Value = '''578.63'
$ws.Range("E5").Value = '  -0.41%  '

# Row 6
$ws.Range("D6").Value = '''141.11'
$ws.Range("E6").Value = '  -1.45%  '

# Row 7
$ws.Range("E7").Value = '  +0.08%  '

# Row 8
$ws.Range("D8").Value = '''0.533'
$ws.Range("E8").Value = '  +0.74%  '

# Row 9
$ws.Range("D9").Value = '2.438.80'
$ws.Range("E9").Value = '  +0.01%  '

# Row 10
$ws.Range("E10").Value = '  +2.61%  '

# Row 11
$ws.Range("E11").Value = '  +1.85%  '

# Row 12
$ws.Range("D12").Value = '''5.17'
$ws.Range("E12").Value = '  -0.71%  '

# Row 13
$ws.Range("D13").Value = '''0.341'
$ws.Range("E13").Value = '  -1.19%  '

# Row 14
$ws.Range("D14").Value = '''25.96'
$ws.Range("E14").Value = '  -1.89%  '

# Row 15
$ws.Range("D15").Value = '2.890.26'
$ws.Range("E15").Value = '  +0.75%  '

# Row 16
$ws.Range("E16").Value = '  -0.45%  '

# Row 17
$ws.Range("D17").Value = '61.773.41'
$ws.Range("E17").Value = '  -0.79%  '

# Row 18
$ws.Range("D18").Value = '2.459.58'
$ws.Range("E18").Value = '  +1.10%  '

# Row 19
$ws.Range("D19").Value = '''10.61'
$ws.Range("E19").Value = '  -3.19%  '

# Row 20
$ws.Range("D20").Value = '''7.25'
$ws.Range("E20").Value = '  +1.99%  '

# Row 21
$ws.Range("D21").Value = '''325.28'
$ws.Range("E21").Value = '  -1.87%  '

# Row 22
$ws.Range("E22").Value = '  -1.01%  '

# Row 23
$ws.Range("D23").Value = '''6.04'
$ws.Range("E23").Value = '  +0.99%  '

# Row 24
$ws.Range("D24").Value = '''1.95'
$ws.Range("E24").Value = '  -0.33%  '

# Row 25
$ws.Range("E25").Value = '  +0.03%  '

# Row 26
$ws.Range("D26").Value = '''64.97'
$ws.Range("E26").Value = '  -1.12%  '

# Row 27
$ws.Range("D27").Value = '''9.13'
$ws.Range("E27").Value = '  -1.48%  '

# Row 28
$ws.Range("D28").Value = '''582.03'
$ws.Range("E28").Value = '  -7.86%  '

# Row 29
$ws.Range("E29").Value = '  +0.13%  '

# Row 30
$ws.Range("E30").Value = '  +0.00%  '

# Row 31
$ws.Range("D31").Value = '0.0₃0922'
$ws.Range("E31").Value = '  -3.16%  '

# Row 32
$ws.Range("D32").Value = '''7.93'
$ws.Range("E32").Value = '  -1.27%  '

# Row 33
$ws.Range("E33").Value = '  -4.73%  '

# Row 34
$ws.Range("E34").Value = '  -0.97%  '

# Row 35
$ws.Range("E35").Value = '  -6.31%  '

# Row 36
$ws.Range("E36").Value = '  +0.01%  '

# Row 37
$ws.Range("D37").Value = '''4.69'
$ws.Range("E37").Value = '  -5.36%  '

# Row 38
$ws.Range("E38").Value = '  -0.67%  '

# Row 39
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '''151.84'
$ws.Range("E39").Value = '  +1.45%  '

# Row 40
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").Value = '''1.40'
$ws.Range("E40").Value = '  -3.21%  '

# Row 41
$ws.Range("D41").Value = '''18.30'
$ws.Range("E41").Value = '  -0.27%  '

# Row 42
$ws.Range("D42").Value = '''5.16'
$ws.Range("E42").Value = '  -1.74%  '

# Row 43
$ws.Range("E43").Value = '  +0.01%  '

# Row 44
$ws.Range("E44").Value = '  -4.15%  '

# Row 45
$ws.Range("D45").Value = '''41.67'
$ws.Range("E45").Value = '  -2.54%  '

# Row 46
$ws.Range("D46").Value = '''2.37'
$ws.Range("E46").Value = '  -4.17%  '

# Row 47
$ws.Range("D47").Value = '0.0₆0291'
$ws.Range("E47").Value = '  +23.59%  '

# Row 48
$ws.Range("D48").Value = '''142.82'
$ws.Range("E48").Value = '  -0.23%  '

# Row 49
$ws.Range("D49").Value = '''3.57'
$ws.Range("E49").Value = '  -2.37%  '

# Row 50
$ws.Range("E50").Value = '  +0.17%  '

# Row 51
$ws.Range("D51").Value = '''19.66'
$ws.Range("E51").Value = '  +0.11%  '
